# POW_Sequencer_Template.xlsx edit script
# - Inserts a new "Configurazione" sheet between "Sequenza" and "Istruzioni"
#   with a source-folder setting table + notes.
# - Rewrites the "Istruzioni" sheet with the expanded guide content
#   (renamed title, configuration/workflow sections, updated examples).

$wb = $excel.ActiveWorkbook
$seq = $wb.Worksheets.Item("Sequenza")

# ---------------------------------------------------------------------------
# 1. Add the "Configurazione" sheet right after "Sequenza"
# ---------------------------------------------------------------------------
$cfg = $wb.Worksheets.Add($null, $seq)
$cfg.Name = "Configurazione"

# Re-resolve "Istruzioni" AFTER the insert: sheet lookups above are
# positional, and inserting a sheet shifts everything after it.
$instr = $wb.Worksheets.Item("Istruzioni")

# Column widths (20 / 50 / 45 characters). ColumnWidth is expressed in the
# same "characters" unit Excel uses for the stored <col width> value, minus
# the standard 5/6 character gridline padding that the COM layer re-adds.
$cfg.Columns.Item(1).ColumnWidth = 20 - 5 / 6
$cfg.Columns.Item(2).ColumnWidth = 50 - 5 / 6
$cfg.Columns.Item(3).ColumnWidth = 45 - 5 / 6

# -- Row 1: header (bold white text on blue fill, bordered, left aligned) --
$cfg.Range("A1").Value = "Impostazione"
$cfg.Range("B1").Value = "Valore"
$cfg.Range("C1").Value = "Descrizione"
$headerRange = $cfg.Range("A1:C1")
$seq.Range("A1").Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats - reuse "Sequenza" header style
$headerRange.HorizontalAlignment = 1  # xlGeneral (not centered, unlike Sequenza)
$wb.Application.CutCopyMode = $false

# -- Row 2: data row (light green fill, bordered) --
$cfg.Range("A2").Value = "Cartella Sorgenti"
$cfg.Range("B2").Value = "default"
$cfg.Range("C2").Value = "Percorso dove salvare i file MDB da Powin-PC2"
$dataRange = $cfg.Range("A2:C2")
$dataRange.Borders.LineStyle = 1
$dataRange.Interior.Color = 14348258  # RGB(0xE2,0xEF,0xDA) light green

# -- Row 4: "NOTE:" bold label --
$cfg.Range("A4").Value = "NOTE:"
$cfg.Range("A4").Font.Bold = $true

# -- Rows 5-7: plain note lines --
$cfg.Range("A5").Value = "- 'default' = usa la cartella 'Sorgenti' nella stessa directory del file Excel"
$cfg.Range("A6").Value = "- Esegui la macro 'SelectSourceFolder' per selezionare una cartella diversa"
$cfg.Range("A7").Value = "- Esegui la macro 'CheckSourceFiles' per verificare lo stato dei file"

$cfg.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Rewrite "Istruzioni" with the expanded guide
# ---------------------------------------------------------------------------
$instr.Cells.Clear()

$instr.Cells.Item(1, 1).Value = "POW PROGRAM SEQUENCER - GUIDA"

$instr.Cells.Item(3, 1).Value = "CONFIGURAZIONE INIZIALE:"
$instr.Cells.Item(4, 1).Value = "1. Crea una cartella 'Sorgenti' nella stessa directory del file Excel"
$instr.Cells.Item(5, 1).Value = "2. Oppure usa la macro 'SelectSourceFolder' per scegliere una cartella diversa"
$instr.Cells.Item(6, 1).Value = "3. Salva i file MDB da Powin-PC2 nella cartella Sorgenti:"
$instr.Cells.Item(7, 1).Value = "   - 30IGNIT.mdb"
$instr.Cells.Item(8, 1).Value = "   - 31NOWELD.mdb"
$instr.Cells.Item(9, 1).Value = "   - 32WELD.mdb"
$instr.Cells.Item(10, 1).Value = "   - 33DWNSLP.mdb"

$instr.Cells.Item(12, 1).Value = "FLUSSO DI LAVORO:"
$instr.Cells.Item(13, 1).Value = "1. Modifica i parametri in Powin-PC2"
$instr.Cells.Item(14, 1).Value = "2. Esporta/Salva il programma modificato come file MDB"
$instr.Cells.Item(15, 1).Value = "3. Copia il file MDB nella cartella Sorgenti"
$instr.Cells.Item(16, 1).Value = "4. Apri questo Excel e imposta la sequenza"
$instr.Cells.Item(17, 1).Value = "5. Esegui 'GenerateMDB' per creare il file combinato"

$instr.Cells.Item(19, 1).Value = "COME USARE:"
$instr.Cells.Item(20, 1).Value = "1. Nel foglio 'Sequenza', inserisci i numeri nella colonna A"
$instr.Cells.Item(21, 1).Value = "2. L'ordine delle righe determina la sequenza di esecuzione"
$instr.Cells.Item(22, 1).Value = "3. Salva il file come .xlsm (con macro)"
$instr.Cells.Item(23, 1).Value = "4. Importa il modulo VBA 'POW_Sequencer_VBA.bas'"
$instr.Cells.Item(24, 1).Value = "5. Esegui la macro 'GenerateMDB'"

$instr.Cells.Item(26, 1).Value = "PROGRAMMI DISPONIBILI:"
$instr.Cells.Item(27, 1).Value = "  30 = IGNIT (Accensione) - 12 funzioni"
$instr.Cells.Item(28, 1).Value = "  31 = NOWELD (No saldatura) - 39 funzioni"
$instr.Cells.Item(29, 1).Value = "  32 = WELD (Saldatura) - 49 funzioni"
$instr.Cells.Item(30, 1).Value = "  33 = DWNSLP (Downslope) - 49 funzioni"

$instr.Cells.Item(32, 1).Value = "MACRO DISPONIBILI:"
$instr.Cells.Item(33, 1).Value = "  GenerateMDB       - Genera il file MDB dalla sequenza"
$instr.Cells.Item(34, 1).Value = "  SelectSourceFolder - Seleziona la cartella sorgenti"
$instr.Cells.Item(35, 1).Value = "  CheckSourceFiles   - Verifica stato dei file sorgente"
$instr.Cells.Item(36, 1).Value = "  ClearSequence      - Pulisce la sequenza"
$instr.Cells.Item(37, 1).Value = "  AddDefaultSequence - Aggiunge sequenza 30-31-32-33"
$instr.Cells.Item(38, 1).Value = "  ShowHelp           - Mostra la guida"

$instr.Cells.Item(40, 1).Value = "ESEMPIO:"
$instr.Cells.Item(41, 1).Value = "Se inserisci: 30, 32, 33"
$instr.Cells.Item(42, 1).Value = "Il programma MDB finale conterra:"
$instr.Cells.Item(43, 1).Value = "  - Prima tutte le funzioni di 30IGNIT (righe 1-11)"
$instr.Cells.Item(44, 1).Value = "  - Poi tutte le funzioni di 32WELD (righe 12-59)"
$instr.Cells.Item(45, 1).Value = "  - Infine tutte le funzioni di 33DWNSLP (righe 60-107)"

# Restore the originally active sheet/selection ("Sequenza" was active before
# this edit; adding/editing sheets above shifts the active tab).
$seq.Activate()
$seq.Range("A1").Select()
